$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$r = $ws.Range("S4")
$r.Value = "x"
$r.Font.Name = "Segoe UI"
$r.Font.Size = 10
$r.Font.Color = 0x00A0A09A
Write-Host "ok"
